# GBDS Credits workbook update | September 2025
# Updates the weekly collection report: new collection-dates, new
# collector entries/amounts for ROUTE 1 & ROUTE 3 "ADD: COLLECTION"
# tables, and refreshes the print area / active view to the ROUTE 3
# block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# ROUTE 1 block (rows 1-11) - "ADD: COLLECTION" table (H:L)
# ---------------------------------------------------------------
# Collection-week date (K1) now recorded as a real date instead of text
$ws.Range("K1").Value = 45878

# Row 3 - collector changed, new S.I. no. / amount
$ws.Range("H3").Value = "JAY TABASA"
$ws.Range("I3").Value = 5603
$ws.Range("J3").Value = 84657

# Row 4 - collector changed, new S.I. no. / amount
$ws.Range("H4").Value = "DELIA VILLARIN"
$ws.Range("I4").Value = 5634
$ws.Range("J4").Value = 66997

# Row 5 - previously blank, now a new collector entry
$ws.Range("H5").Value = "SOUND CHECK"
$ws.Range("I5").Value = 5793
$ws.Range("J5").Value = 22825

# ---------------------------------------------------------------
# ROUTE 3 block (rows 25-35) - "ADD: COLLECTION" table (H:L)
# ---------------------------------------------------------------
# Collection-week date (K25) moved forward to the new week
$ws.Range("K25").Value = 45909

# Row 27 - collector changed, new S.I. no. / amount
$ws.Range("H27").Value = "YAKEN STORE"
$ws.Range("I27").Value = 5459
$ws.Range("J27").Value = 229115

# Row 28 - collector changed, S.I. no. cleared, new amount
$ws.Range("H28").Value = "MERIAM APDUHAN"
$ws.Range("I28").ClearContents()
$ws.Range("J28").Value = 50905

# ---------------------------------------------------------------
# Workbook-level: print area now targets the ROUTE 3 block
# ---------------------------------------------------------------
$ws.PageSetup.PrintArea = "`$H`$25:`$M`$36"

# ---------------------------------------------------------------
# View state: zoomed out and scrolled/selected onto the ROUTE 3 block
# ---------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("J26:L26").Select() | Out-Null
